# Weekly data refresh: insert the newest week's two records (Primera /
# Segunda quality grades for Cilantro sold as "atado" in Región de Ñuble)
# at their sorted position (row 235), pushing the rest of the historical
# rows (old rows 235-337) down by two rows to 237-339.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 235 (shifts old row 235.. down to 237..)
$ws.Rows.Item(235).EntireRow.Insert()
$ws.Rows.Item(235).EntireRow.Insert()

# New row 235: "Primera" grade record for the new week
$ws.Range("A235").Value = 11
$ws.Range("B235").Value = "Vega Monumental Concepción"
$ws.Range("C235").Value = "Bíobío"
$ws.Range("D235").Value = 45141
$ws.Range("E235").Value = 8
$ws.Range("F235").Value = 100112040
$ws.Range("G235").Value = "Cilantro"
$ws.Range("H235").Value = "Sin especificar"
$ws.Range("I235").Value = "Primera"
$ws.Range("J235").Value = 200
$ws.Range("K235").Value = 600
$ws.Range("L235").Value = 700
$ws.Range("M235").Value = 650
$ws.Range("N235").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O235").Value = "Región de Ñuble"
$ws.Range("P235").Value = 650
$ws.Range("Q235").Value = 1
$ws.Range("R235").Value = "Hortaliza"

# New row 236: "Segunda" grade record for the new week
$ws.Range("A236").Value = 11
$ws.Range("B236").Value = "Vega Monumental Concepción"
$ws.Range("C236").Value = "Bíobío"
$ws.Range("D236").Value = 45141
$ws.Range("E236").Value = 8
$ws.Range("F236").Value = 100112040
$ws.Range("G236").Value = "Cilantro"
$ws.Range("H236").Value = "Sin especificar"
$ws.Range("I236").Value = "Segunda"
$ws.Range("J236").Value = 100
$ws.Range("K236").Value = 500
$ws.Range("L236").Value = 500
$ws.Range("M236").Value = 500
$ws.Range("N236").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O236").Value = "Región de Ñuble"
$ws.Range("P236").Value = 500
$ws.Range("Q236").Value = 1
$ws.Range("R236").Value = "Hortaliza"
